# Apply the two changes from the commit:
#   1. Slide 6's table switches to a different table style (GUID).
#   2. The deck's theme colour scheme is swapped from the "Integral" palette
#      to the standard "Office" palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{7AE27441-729D-4819-BC67-12F59084C8E3}")
    }
}

# --- 2. Theme colour scheme: Integral -> Office ---------------------------------
# RGB() in VBA/PowerPoint packs colour bytes as 0x00BBGGRR, i.e. r + g*256 + b*65536.
$officeColors = @(
    0,         # 1  dk1      000000
    16777215,  # 2  lt1      FFFFFF
    6968388,   # 3  dk2      44546A
    15132391,  # 4  lt2      E7E6E6
    13998939,  # 5  accent1  5B9BD5
    3243501,   # 6  accent2  ED7D31
    10855845,  # 7  accent3  A5A5A5
    49407,     # 8  accent4  FFC000
    12874308,  # 9  accent5  4472C4
    4697456,   # 10 accent6  70AD47
    12673797,  # 11 hlink    0563C1
    7491477    # 12 folHlink 954F72
)

$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$theme = $master.Theme
$tcs = $theme.ThemeColorScheme

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
